$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '70.323.41'
$ws.Range('E2').Value = '  -2.63%  '

# Row 3
$ws.Range('D3').Value = '2.523.64'
$ws.Range('E3').Value = '  -4.84%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '575.97'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.54%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '169.56'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.96%  '

# Row 7
$ws.Range('E7').Value = '  +0.12%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.511'
$ws.Range('D8').Style = "Normal"

# Row 9
$ws.Range('D9').Value = '2.522.18'
$ws.Range('E9').Value = '  -4.83%  '

# Row 10
$ws.Range('E10').Value = '  -5.62%  '

# Row 11
$ws.Range('E11').Value = '  -0.72%  '

# Row 12
$ws.Range('E12').Value = '  -3.58%  '

# Row 13
$ws.Range('E13').Value = '  -3.18%  '

# Row 14
$ws.Range('D14').Value = '2.986.15'
$ws.Range('E14').Value = '  -4.79%  '

# Row 15
$ws.Range('D15').Value = '70.153.73'
$ws.Range('E15').Value = '  -2.72%  '

# Row 16
$ws.Range('E16').Value = '  -3.40%  '

# Row 17
$ws.Range('E17').Value = '  -4.41%  '

# Row 18
$ws.Range('D18').Value = '2.520.59'
$ws.Range('E18').Value = '  -4.73%  '

# Row 19
$ws.Range('E19').Value = '  -5.64%  '

# Row 20
$ws.Range('E20').Value = '  -6.49%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '355.48'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.99%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.96'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -5.30%  '

# Row 23
$ws.Range('E23').Value = '  -2.72%  '

# Row 24
$ws.Range('E24').Value = '  +0.02%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '69.08'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -4.12%  '

# Row 26
$ws.Range('E26').Value = '  -4.98%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.19'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -5.47%  '

# Row 28
$ws.Range('E28').Value = '  -4.93%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.995'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.54%  '

# Row 30
$ws.Range('D30').Value = '0.0₃0913'
$ws.Range('E30').Value = '  -5.69%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.86'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.81%  '

# Row 32
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.32'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.19%  '

# Row 33
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '483.65'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.95%  '

# Row 34
$ws.Range('E34').Value = '  -3.19%  '

# Row 35
$ws.Range('E35').Value = '  -0.04%  '

# Row 36
$ws.Range('E36').Value = '  +5.22%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '155.54'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -4.53%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '18.89'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.19%  '

# Row 41
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.321'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.15%  '

# Row 42
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.65'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -6.61%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '4.75'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -4.63%  '

# Row 44
$ws.Range('E44').Value = '  -11.43%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.39'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -7.35%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '38.29'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.94%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '142.69'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -8.14%  '

# Row 48
$ws.Range('E48').Value = '  -5.44%  '

# Row 49
$ws.Range('E49').Value = '  -5.27%  '

# Row 50
$ws.Range('E50').Value = '  -5.98%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.599'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.69%  '

